$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows whose op_type text changed
$ws.Range("B10").Value = 'data queries,geometry measurement,data queries,vegetation 애매함. osm만으로 불가.'
$ws.Range("B13").Value = 'topography,data model conversion,overlay analysis,data queries'
$ws.Range("B14").Value = 'data queries,buffer,buffer,overlay analysis,overlay analysis,two set of geometries'
$ws.Range("B15").Value = 'data queries,overlay analysis'
$ws.Range("B16").Value = 'data editing,data queries,network analysis,classification,data queries,overlay analysis,data queries,overlay analysis,data queries,overlay analysis'
$ws.Range("B17").Value = 'data queries,buffer,overlay analysis,정확히 카메라가 무엇을 뜻하는가'
$ws.Range("B18").Value = 'data queries,buffer,overlay analysis,각자 buffer를 union 해야하네'
$ws.Range("B19").Value = 'data queries,buffer,overlay analysis,matter of scale. 도로는 line인가'
$ws.Range("B20").Value = 'data queries,buffer,overlay analysis,major가 무엇인지'
$ws.Range("B21").Value = 'data queries,network analysis,data queries,network analysis,classification,data queries,overlay analysis,어디로 부터 가장 가까운 소방서'
$ws.Range("B22").Value = 'data queries,network analysis,classification,data queries,overlay analysis,data queries,overlay analysis,data queries,overlay analysis'
$ws.Range("B23").Value = 'data queries,buffer,overlay analysis,shop이 너무나 많다'
$ws.Range("B24").Value = 'data editing,data queries,buffer,overlay analysis'
$ws.Range("B25").Value = 'data editing,data queries,buffer,overlay analysis,what area는 그냥 boundary만 얘기하는 것인가? 아니면 다른 attribute도 clip하라는 것인가. 일단 clip. urban tag를 그렇게 사용하지도 않는다'
$ws.Range("B26").Value = 'topography,classification,data queries,data model conversion,overlay analysis'
$ws.Range("B27").Value = 'data queries,geometry measurement,data editing,data queries'
$ws.Range("B28").Value = 'geostatistics  ,classification,data queries,data model conversion,overlay analysis'
$ws.Range("B29").Value = 'data queries,network analysis,classification,data queries,overlay analysis,data queries,overlay analysis'
$ws.Range("B30").Value = 'data queries,buffer,overlay analysis,data queries'
$ws.Range("B31").Value = 'network analysis,data queries,buffer,overlay analysis,data queries'
$ws.Range("B32").Value = 'data queries,overlay analysis,data queries,osm urban이 있긴 하지만 거의 안쓴다'
$ws.Range("B33").Value = 'overlay analysis,topography'
$ws.Range("B34").Value = 'data queries,overlay analysis,data editing,data queries'
$ws.Range("B35").Value = 'data editing,data queries,overlay analysis,data editing,data queries'
$ws.Range("B36").Value = 'data editing,overlay analysis,data editing,data queries'
$ws.Range("B37").Value = 'data editing,data queries,data editing'
$ws.Range("B38").Value = 'data queries,generalization,geostatistics  ,https://pro.arcgis.com/en/pro-app/latest/tool-reference/spatial-statistics/h-how-central-feature-spatial-statistics-works.htm'
$ws.Range("B39").Value = 'data editing,data queries,generalization,geostatistics  '
$ws.Range("B40").Value = 'data queries,generalization,geostatistics  '
$ws.Range("B41").Value = 'data editing,overlay analysis,data queries'
$ws.Range("B42").Value = 'data editing,data queries,geostatistics  '

# Add new rows 43-56 with A (index), B (op_type text), C (op_id) and copy formatting for column A from row 2
$ws.Range("A43").Value = 41
$ws.Range("B43").Value = 'data queries,geostatistics  '
$ws.Range("C43").Value = 41
$ws.Range("A44").Value = 42
$ws.Range("B44").Value = 'data editing,data queries,overlay analysis,data editing'
$ws.Range("C44").Value = 42
$ws.Range("A45").Value = 43
$ws.Range("B45").Value = 'data editing,data queries,network analysis,data queries'
$ws.Range("C45").Value = 43
$ws.Range("A46").Value = 44
$ws.Range("B46").Value = 'data queries,overlay analysis,geostatistics  '
$ws.Range("C46").Value = 44
$ws.Range("A47").Value = 45
$ws.Range("B47").Value = 'data queries,network analysis'
$ws.Range("C47").Value = 45
$ws.Range("A48").Value = 46
$ws.Range("B48").Value = 'data editing,buffer,overlay analysis,data editing,data queries'
$ws.Range("C48").Value = 46
$ws.Range("A49").Value = 47
$ws.Range("B49").Value = 'data queries,geometry measurement,data queries'
$ws.Range("C49").Value = 47
$ws.Range("A50").Value = 48
$ws.Range("B50").Value = 'data queries,overlay analysis,data queries,어디선 polygon, 어디선 linestring'
$ws.Range("C50").Value = 48
$ws.Range("A51").Value = 49
$ws.Range("B51").Value = 'data queries,geometry measurement,data queries,https://wiki.openstreetmap.org/wiki/sparql_examples'
$ws.Range("C51").Value = 49
$ws.Range("A52").Value = 50
$ws.Range("B52").Value = 'data queries,buffer,overlay analysis,data queries,geometry measurement,data queries'
$ws.Range("C52").Value = 50
$ws.Range("A53").Value = 51
$ws.Range("B53").Value = 'data queries,buffer,buffer,buffer,overlay analysis,overlay analysis,overlay analysis,geometry measurement,data queries,geometry measurement,data queries,geometry measurement,data queries'
$ws.Range("C53").Value = 51
$ws.Range("A54").Value = 52
$ws.Range("B54").Value = 'data queries,geometry measurement,data queries,buffer,overlay analysis,data queries'
$ws.Range("C54").Value = 52
$ws.Range("A55").Value = 53
$ws.Range("B55").Value = 'data queries,network analysis,data queries,network analysis,classification,data queries,overlay analysis'
$ws.Range("C55").Value = 53
$ws.Range("A56").Value = 54
$ws.Range("B56").Value = 'geocoding,data queries,network analysis,data queries,멘탈 헬스를 누가 담당할건데…'
$ws.Range("C56").Value = 54

# Copy the style (border/font/alignment) from an existing formatted A cell to the new A cells
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A43:A56").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("A1").Select() | Out-Null